# Commit for safety before Interfaces
#
# Adds a new diary entry (row 23) for 13.11.18 - "-Mouse Button Simulation",
# and corrects the end time of the 12.11.18 entry (row 22, column C) from
# 18:00 to 18:45, which ripples through the "time" formula column (D) and
# the running total in D26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the end-time of the 12.11.18 row (row 22): 18:00 -> 18:45
$ws.Range("C22").Value = 0.78125

# Fill in the new diary entry on row 23: 13.11.18, 18:30 - 21:30, Mouse Button Simulation
$ws.Range("A23").Value = "13.11.18"
$ws.Range("B23").Value = 0.77083333333333337
$ws.Range("C23").Value = 0.89583333333333337
$ws.Range("E23").Value = "-Mouse Button Simulation"

# Move the active selection to C24, matching the post-edit cursor position
$ws.Range("C24").Select()
